$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.085355997085571
$ws.Range("B1").Value = 2.970558881759644
$ws.Range("C1").Value = 5.772488117218018
$ws.Range("D1").Value = 2.748442411422729
$ws.Range("E1").Value = 1.182768583297729
